# edit.ps1 - PowerPoint COM-interop script
#
# Reproduces two changes from the source diff:
#
#  1. On slide 16, the table's table-style is switched from the
#     deck's custom "Table_0" style ({1D09F7B5-11E3-4B42-B070-97CC4070FA6F})
#     to the built-in style {B476317D-7368-44DE-A798-5058CBBAB25E}.
#
#  2. The presentation's active theme colour scheme is changed from the
#     "Integral" palette to the standard Office "Office Theme" palette
#     (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 16 -------------------------------
$slide = $p.Slides.Item(16)
$tableShape = $slide.Shapes.Item(3)
$tableShape.Table.ApplyStyle("{B476317D-7368-44DE-A798-5058CBBAB25E}")

# --- 2. Swap the live theme colours: Integral -> Office Theme --------
# Helper: turn an "RRGGBB" hex string into the BGR-packed integer that
# the ThemeColor.RGB property (and VBA's RGB()) expects.
function ConvertTo-ComRgb($hex) {
    $v = [Convert]::ToInt32($hex, 16)
    $r = ($v -band 0xFF0000) -shr 16
    $g = ($v -band 0x00FF00) -shr 8
    $b = ($v -band 0x0000FF)
    return ($b * 65536) + ($g * 256) + $r
}

# Index order matches ThemeColorScheme.Colors(1..12):
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5, accent6, hlink, folHlink
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$themeColors = $slide.ThemeColorScheme

for ($idx = 1; $idx -le 12; $idx++) {
    $hexVal = $officeThemeColors[$idx - 1]
    $themeColors.Colors($idx).RGB = ConvertTo-ComRgb $hexVal
}
